# "Generate Report for Handoff"
# Adds a new localization-status row (file 85aa701b-b5f6-4eb6-ab99-26ea18c264ab...md,
# which is "Ready for handoff") to the Overview / zh-cn / de-de tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Shared literal values used across sheets
# ---------------------------------------------------------------------------
$fileName        = "85aa701b-b5f6-4eb6-ab99-26ea18c264aboooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$pathAndName     = "e2e\85aa701b-b5f6-4eb6-ab99-26ea18c264aboooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$extension       = ".md"
$statusText      = "Ready for handoff"
$hoDate          = "2016-09-06 16:54:59"

$zhXlf           = "85aa701b-b5f6-4eb6-ab99-26ea18c264abooooooooooooooooooooooooooooooooooooooo.8567f090ea009d28fee48faaf22ef6a1e108f984.zh-cn.xlf"
$zhHandoffDate   = "2016-09-06 16:54:47"
$deXlf           = "85aa701b-b5f6-4eb6-ab99-26ea18c264abooooooooooooooooooooooooooooooooooooooo.8567f090ea009d28fee48faaf22ef6a1e108f984.de-de.xlf"
$deHandoffDate   = $hoDate

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a773427fe68b0647af1fa1e9c27b40d05f41574/e2e/" + $fileName

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $fileName
$wsOverview.Range("B3").Value = $pathAndName
$wsOverview.Range("C3").Value = $extension
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Range("G3").Value = $hoDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkTarget, "", "", $pathAndName) | Out-Null

# Status column values got longer ("Ready for handoff" > "In Translation"), so
# the sheet got re-auto-fit - widen E/F to match.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333336

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $fileName
$wsZh.Range("B3").Value = $extension
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $hyperlinkTarget, "", "", $fileName) | Out-Null

$wsZh.Columns.Item(3).ColumnWidth = 16.333333333333336

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $fileName
$wsDe.Range("B3").Value = $extension
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $hyperlinkTarget, "", "", $fileName) | Out-Null

$wsDe.Columns.Item(3).ColumnWidth = 16.333333333333336
